# Updated cryptos list with GitHub Actions
# Refreshes price / volume(1h) figures in columns D/E, and re-ranks a
# few coins (Aave <-> BabyDogeCoin swap rows 46/47; Mantle replaced by
# RenderToken in row 51) to mirror a fresh scrape of coinranking.com.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "26.200.22"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  -1.64%  "

# Row 3
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.660.66"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  -1.52%  "

# Row 4
$ws.Range("E4").Value = "  +0.37%  "

# Row 5
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "218.04"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +0.45%  "

# Row 6
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "0.5215"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -2.36%  "

# Row 7
$ws.Range("E7").Value = "  +0.36%  "

# Row 8
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.2664"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  -0.54%  "

# Row 9
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.06309"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  -1.99%  "

# Row 10
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "21.02"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  -2.98%  "

# Row 11
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.07717"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -1.16%  "

# Row 12
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "1.650.30"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  -2.25%  "

# Row 13
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "4.425"
$cell.Style = "Normal"

# Row 14
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "1.885.26"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -1.66%  "

# Row 15
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.5461"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  -2.67%  "

# Row 16
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "0.0₅8203"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  -2.84%  "

# Row 17
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "64.74"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  -1.95%  "

# Row 18
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "26.232.86"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  -1.66%  "

# Row 19
$ws.Range("E19").Value = "  +0.39%  "

# Row 20
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "4.659"
$cell.Style = "Normal"

# Row 21
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "192.35"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -1.53%  "

# Row 22
$ws.Range("E22").Value = "  -2.66%  "

# Row 23
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "6.076"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -4.70%  "

# Row 24
$ws.Range("E24").Value = "  +0.52%  "

# Row 25
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "138.64"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -3.97%  "

# Row 26
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "0.1238"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  -3.44%  "

# Row 27
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "7.219"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -3.55%  "

# Row 28
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "16.14"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -1.23%  "

# Row 29
$ws.Range("E29").Value = "  -1.37%  "

# Row 30
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "0.06008"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  -2.39%  "

# Row 31
$ws.Range("E31").Value = "  +0.17%  "

# Row 32
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "3.605"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -0.16%  "

# Row 33
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "3.315"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -4.47%  "

# Row 34
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "1.637"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -3.79%  "

# Row 35
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "0.9779"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -3.36%  "

# Row 37
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "2.782"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -0.60%  "

# Row 38
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.5878"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  +2.21%  "

# Row 39
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "0.01587"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -3.63%  "

# Row 40
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "5.942"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -1.59%  "

# Row 41
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.8638"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  +0.03%  "

# Row 43
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "1.031.89"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  -3.66%  "

# Row 44
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "99.58"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -0.86%  "

# Row 45
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "1.801.07"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -1.96%  "

# Row 46
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "56.98"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -0.49%  "

# Row 47
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "0.0₈107"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -1.90%  "

# Row 48
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "1.010"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  +0.67%  "

# Row 49
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "8.095"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  -1.15%  "

# Row 50
$ws.Range("E50").Value = "  -0.60%  "

# Row 51
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "1.463"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +0.30%  "
